$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2594.817
$ws.Range("I15").Value = 2594.817
$ws.Range("K15").Value = 7784.451
$ws.Range("M15").Value = -7615.451
# Row 17
$ws.Range("H17").Value = 2020.2122
$ws.Range("J17").Value = 1812.3103
$ws.Range("L17").Value = 5436.9309
$ws.Range("N17").Value = -5772.9309
# Row 112
$ws.Range("H112").Value = 2379.5833
$ws.Range("J112").Value = 2379.5833
$ws.Range("L112").Value = 7138.749899999999
$ws.Range("N112").Value = -9354.749899999999
# Row 132
$ws.Range("H132").Value = 3573.9048
$ws.Range("I132").Value = 3655.5615
$ws.Range("K132").Value = 10966.6845
$ws.Range("M132").Value = -8436.684499999999
# Row 138
$ws.Range("H138").Value = 3113.8713
$ws.Range("I138").Value = 2726.875
$ws.Range("J138").Value = 3315.7827
$ws.Range("K138").Value = 8180.625
$ws.Range("L138").Value = 9947.348100000001
$ws.Range("M138").Value = -3040.625
$ws.Range("N138").Value = -20227.3481

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 63202.145
$ws.Range("I45").Value = 142106.67
$ws.Range("J45").Value = 4023.75
$ws.Range("K45").Value = 142106.67
$ws.Range("L45").Value = 4023.75
$ws.Range("M45").Value = -141729.67
$ws.Range("N45").Value = -4777.75
# Row 61
$ws.Range("H61").Value = 4910.737
$ws.Range("I61").Value = 2661.3635
$ws.Range("K61").Value = 2661.3635
$ws.Range("M61").Value = -2449.3635
# Row 130
$ws.Range("H130").Value = 38196.4
$ws.Range("J130").Value = 38196.4
$ws.Range("L130").Value = 38196.4
$ws.Range("N130").Value = -48236.4
# Row 132
$ws.Range("H132").Value = 5210903.5
$ws.Range("I132").Value = 2249.5745
$ws.Range("J132").Value = 19611300
$ws.Range("K132").Value = 6748.7235
$ws.Range("L132").Value = 58833900
$ws.Range("M132").Value = -4218.7235
$ws.Range("N132").Value = -58838960
# Row 136
$ws.Range("H136").Value = 4910.737
$ws.Range("I136").Value = 2661.3635
$ws.Range("K136").Value = 7984.0905
$ws.Range("M136").Value = -5434.0905

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 77300
$ws.Range("J2").Value = 77300
$ws.Range("L2").Value = 77300
$ws.Range("N2").Value = -77526
# Row 86
$ws.Range("H86").Value = 2856.5
$ws.Range("I86").Value = 2462.25
$ws.Range("J86").Value = 3382.1667
$ws.Range("K86").Value = 2462.25
$ws.Range("L86").Value = 3382.1667
$ws.Range("M86").Value = -1339.25
$ws.Range("N86").Value = -5628.1667
# Row 89
$ws.Range("H89").Value = 2856.5
$ws.Range("I89").Value = 2462.25
$ws.Range("J89").Value = 3382.1667
$ws.Range("K89").Value = 12311.25
$ws.Range("L89").Value = 16910.8335
$ws.Range("M89").Value = -6695.25
$ws.Range("N89").Value = -28142.8335
# Row 124
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1528.4615
$ws.Range("I22").Value = 1622.625
$ws.Range("K22").Value = 1622.625
$ws.Range("M22").Value = -1272.625
# Row 31
$ws.Range("H31").Value = 2320527.5
$ws.Range("I31").Value = 5327.96
$ws.Range("J31").Value = 4316389
$ws.Range("K31").Value = 5327.96
$ws.Range("L31").Value = 4316389
$ws.Range("M31").Value = -5032.96
$ws.Range("N31").Value = -4316979
# Row 34
$ws.Range("H34").Value = 2320527.5
$ws.Range("I34").Value = 5327.96
$ws.Range("J34").Value = 4316389
$ws.Range("K34").Value = 5327.96
$ws.Range("L34").Value = 4316389
$ws.Range("M34").Value = -5125.96
$ws.Range("N34").Value = -4316793
# Row 36
$ws.Range("H36").Value = 12496.5
$ws.Range("J36").Value = 12496.5
$ws.Range("L36").Value = 12496.5
$ws.Range("N36").Value = -13272.5
# Row 40
$ws.Range("H40").Value = 12496.5
$ws.Range("J40").Value = 12496.5
$ws.Range("L40").Value = 12496.5
$ws.Range("N40").Value = -12816.5
# Row 58
$ws.Range("H58").Value = 1750.9117
$ws.Range("I58").Value = 1147.3462
$ws.Range("J58").Value = 3712.5
$ws.Range("K58").Value = 1147.3462
$ws.Range("L58").Value = 3712.5
$ws.Range("M58").Value = -944.3462
$ws.Range("N58").Value = -4118.5
# Row 86
$ws.Range("H86").Value = 9889.5
$ws.Range("I86").Value = 9889.5
$ws.Range("K86").Value = 9889.5
$ws.Range("M86").Value = -8766.5
# Row 89
$ws.Range("H89").Value = 9889.5
$ws.Range("I89").Value = 9889.5
$ws.Range("K89").Value = 49447.5
$ws.Range("M89").Value = -43831.5
# Row 132
$ws.Range("H132").Value = 5210757.5
$ws.Range("I132").Value = 2022.58
$ws.Range("K132").Value = 6067.74
$ws.Range("M132").Value = -3537.74
# Row 134
$ws.Range("H134").Value = 2557.439
$ws.Range("I134").Value = 2513.0881
$ws.Range("J134").Value = 2772.8572
$ws.Range("K134").Value = 7539.2643
$ws.Range("L134").Value = 8318.571599999999
$ws.Range("M134").Value = -5004.2643
$ws.Range("N134").Value = -13388.5716
# Row 136
$ws.Range("H136").Value = 1750.9117
$ws.Range("I136").Value = 1147.3462
$ws.Range("J136").Value = 3712.5
$ws.Range("K136").Value = 3442.0386
$ws.Range("L136").Value = 11137.5
$ws.Range("M136").Value = -892.0385999999999
$ws.Range("N136").Value = -16237.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 2746.25
$ws.Range("J55").Value = 2746.25
$ws.Range("L55").Value = 8238.75
$ws.Range("N55").Value = -8592.75
# Row 62
$ws.Range("H62").Value = 7904.6665
$ws.Range("J62").Value = 7904.6665
$ws.Range("L62").Value = 23713.9995
$ws.Range("N62").Value = -25085.9995
# Row 65
$ws.Range("H65").Value = 7904.6665
$ws.Range("J65").Value = 7904.6665
$ws.Range("L65").Value = 71141.9985
$ws.Range("N65").Value = -78005.9985
# Row 112
$ws.Range("H112").Value = 2893.3333
$ws.Range("I112").Value = 2893.3333
$ws.Range("K112").Value = 8679.999899999999
$ws.Range("M112").Value = -7571.999899999999
# Row 122
$ws.Range("H122").Value = 640.75
$ws.Range("J122").Value = 917.9167
$ws.Range("L122").Value = 8261.2503
$ws.Range("N122").Value = -13161.2503
# Row 128
$ws.Range("H128").Value = 268605.8
$ws.Range("I128").Value = 268605.8
$ws.Range("K128").Value = 805817.3999999999
$ws.Range("M128").Value = -800837.3999999999
# Row 132
$ws.Range("H132").Value = 2825
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2825
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 25425
$ws.Range("N132").Value = -30485
$ws.Range("M132").ClearContents()
# Row 133
$ws.Range("H133").Value = 4870.1113
$ws.Range("I133").Value = 4449
$ws.Range("J133").Value = 4990.4287
$ws.Range("K133").Value = 13347
$ws.Range("L133").Value = 14971.2861
$ws.Range("M133").Value = -8287
$ws.Range("N133").Value = -25091.2861
# Row 140
$ws.Range("H140").Value = 9001.947
$ws.Range("I140").Value = 9760.75
$ws.Range("K140").Value = 29282.25
$ws.Range("M140").Value = -24102.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Range("H32").Value = 34146.332
$ws.Range("I32").Value = 33054.332
$ws.Range("J32").Value = 34692.332
$ws.Range("K32").Value = 33054.332
$ws.Range("L32").Value = 34692.332
$ws.Range("M32").Value = -32758.332
$ws.Range("N32").Value = -35284.332
# Row 52
$ws.Range("H52").Value = 49999.668
$ws.Range("J52").Value = 49999.668
$ws.Range("L52").Value = 49999.668
$ws.Range("N52").Value = -50517.668
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 120
$ws.Range("H120").Value = 155999
$ws.Range("J120").Value = 155999
$ws.Range("L120").Value = 155999
$ws.Range("N120").Value = -165675

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 6133.7144
$ws.Range("I132").Value = 2747.5833
$ws.Range("K132").Value = 8242.749899999999
$ws.Range("M132").Value = -5712.749899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 18302.75
$ws.Range("J41").Value = 18302.75
$ws.Range("L41").Value = 18302.75
$ws.Range("N41").Value = -19082.75
# Row 118
$ws.Range("H118").Value = 100000
$ws.Range("J118").Value = 100000
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314
# Row 122
$ws.Range("H122").Value = 31251938
$ws.Range("I122").Value = 2401.5
$ws.Range("K122").Value = 7204.5
$ws.Range("M122").Value = -4754.5
# Row 126
$ws.Range("H126").Value = 1259.8
$ws.Range("I126").Value = 1233.3334
$ws.Range("K126").Value = 3700.0002
$ws.Range("M126").Value = -1230.0002
